$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "대체교과목번호" (alternate subject code) placeholder values out of
# F2:F15 - these cells keep their style (s="1") but lose their shared-string
# content, becoming blank cells.
$ws.Range("F2:F15").ClearContents()

# Move / record the active selection at F15 (previously F18).
$ws.Range("F15").Select() | Out-Null
